$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows for RM 232 and SC 92 (entire rows removed from dataset)
$ws.Rows.Item(26).Delete()
$ws.Rows.Item(27).Delete()

# Set all data values for rows 2-33 (A:F) to match target state
$ws.Cells.Item(2, 1).Value = "RM 2"
$ws.Cells.Item(2, 2).Value = -19.7
$ws.Cells.Item(2, 3).Value = 14.9
$ws.Cells.Item(2, 4).Value = -13.5
$ws.Cells.Item(2, 5).Value = -7.2
$ws.Cells.Item(2, 6).Value = 18.03
$ws.Cells.Item(3, 1).Value = "RM 8"
$ws.Cells.Item(3, 2).Value = -19.7
$ws.Cells.Item(3, 3).Value = 11.2
$ws.Cells.Item(3, 4).Value = ""
$ws.Cells.Item(3, 5).Value = ""
$ws.Cells.Item(3, 6).Value = 17.64
$ws.Cells.Item(4, 1).Value = "RM 9"
$ws.Cells.Item(4, 2).Value = -18.7
$ws.Cells.Item(4, 3).Value = ""
$ws.Cells.Item(4, 4).Value = -15.4
$ws.Cells.Item(4, 5).Value = -6.4
$ws.Cells.Item(4, 6).Value = ""
$ws.Cells.Item(5, 1).Value = "RM 14"
$ws.Cells.Item(5, 2).Value = -19.5
$ws.Cells.Item(5, 3).Value = 12.3
$ws.Cells.Item(5, 4).Value = ""
$ws.Cells.Item(5, 5).Value = -5
$ws.Cells.Item(5, 6).Value = 17.66
$ws.Cells.Item(6, 1).Value = "RM 21"
$ws.Cells.Item(6, 2).Value = -19.8
$ws.Cells.Item(6, 3).Value = ""
$ws.Cells.Item(6, 4).Value = -14.2
$ws.Cells.Item(6, 5).Value = -5.7
$ws.Cells.Item(6, 6).Value = 16.43
$ws.Cells.Item(7, 1).Value = "RM 32"
$ws.Cells.Item(7, 2).Value = -19.5
$ws.Cells.Item(7, 3).Value = 15
$ws.Cells.Item(7, 4).Value = -13.8
$ws.Cells.Item(7, 5).Value = -7.1
$ws.Cells.Item(7, 6).Value = 17.24
$ws.Cells.Item(8, 1).Value = "RM 38"
$ws.Cells.Item(8, 2).Value = -19.9
$ws.Cells.Item(8, 3).Value = 15.5
$ws.Cells.Item(8, 4).Value = ""
$ws.Cells.Item(8, 5).Value = ""
$ws.Cells.Item(8, 6).Value = ""
$ws.Cells.Item(9, 1).Value = "RM 42"
$ws.Cells.Item(9, 2).Value = -20.6
$ws.Cells.Item(9, 3).Value = 10.5
$ws.Cells.Item(9, 4).Value = ""
$ws.Cells.Item(9, 5).Value = ""
$ws.Cells.Item(9, 6).Value = 17.26
$ws.Cells.Item(10, 1).Value = "RM 52 a"
$ws.Cells.Item(10, 2).Value = -19.8
$ws.Cells.Item(10, 3).Value = 11.5
$ws.Cells.Item(10, 4).Value = -14.7
$ws.Cells.Item(10, 5).Value = -6.1
$ws.Cells.Item(10, 6).Value = ""
$ws.Cells.Item(11, 1).Value = "RM 58"
$ws.Cells.Item(11, 2).Value = ""
$ws.Cells.Item(11, 3).Value = 11.4
$ws.Cells.Item(11, 4).Value = -15.5
$ws.Cells.Item(11, 5).Value = -7.9
$ws.Cells.Item(11, 6).Value = 17.65
$ws.Cells.Item(12, 1).Value = "RM 81"
$ws.Cells.Item(12, 2).Value = ""
$ws.Cells.Item(12, 3).Value = 12.5
$ws.Cells.Item(12, 4).Value = -14.1
$ws.Cells.Item(12, 5).Value = -5.3
$ws.Cells.Item(12, 6).Value = ""
$ws.Cells.Item(13, 1).Value = "RM 88"
$ws.Cells.Item(13, 2).Value = -19.9
$ws.Cells.Item(13, 3).Value = 12.5
$ws.Cells.Item(13, 4).Value = -13.9
$ws.Cells.Item(13, 5).Value = -5.3
$ws.Cells.Item(13, 6).Value = ""
$ws.Cells.Item(14, 1).Value = "RM 90"
$ws.Cells.Item(14, 2).Value = -19.6
$ws.Cells.Item(14, 3).Value = ""
$ws.Cells.Item(14, 4).Value = -13.1
$ws.Cells.Item(14, 5).Value = -5.4
$ws.Cells.Item(14, 6).Value = 17.76
$ws.Cells.Item(15, 1).Value = "RM 95"
$ws.Cells.Item(15, 2).Value = ""
$ws.Cells.Item(15, 3).Value = 12.5
$ws.Cells.Item(15, 4).Value = -15.2
$ws.Cells.Item(15, 5).Value = -8.4
$ws.Cells.Item(15, 6).Value = 16.2
$ws.Cells.Item(16, 1).Value = "RM 103"
$ws.Cells.Item(16, 2).Value = -19.5
$ws.Cells.Item(16, 3).Value = 13.5
$ws.Cells.Item(16, 4).Value = -15.3
$ws.Cells.Item(16, 5).Value = -6.9
$ws.Cells.Item(16, 6).Value = 17.34
$ws.Cells.Item(17, 1).Value = "RM 116"
$ws.Cells.Item(17, 2).Value = -19.4
$ws.Cells.Item(17, 3).Value = 11.2
$ws.Cells.Item(17, 4).Value = -14.7
$ws.Cells.Item(17, 5).Value = -7.3
$ws.Cells.Item(17, 6).Value = 17.78
$ws.Cells.Item(18, 1).Value = "RM 120"
$ws.Cells.Item(18, 2).Value = -19.6
$ws.Cells.Item(18, 3).Value = 11.5
$ws.Cells.Item(18, 4).Value = -15.2
$ws.Cells.Item(18, 5).Value = -8.5
$ws.Cells.Item(18, 6).Value = ""
$ws.Cells.Item(19, 1).Value = "RM 125"
$ws.Cells.Item(19, 2).Value = -20.6
$ws.Cells.Item(19, 3).Value = 13.2
$ws.Cells.Item(19, 4).Value = ""
$ws.Cells.Item(19, 5).Value = ""
$ws.Cells.Item(19, 6).Value = 17.81
$ws.Cells.Item(20, 1).Value = "RM 134"
$ws.Cells.Item(20, 2).Value = -19
$ws.Cells.Item(20, 3).Value = 12.5
$ws.Cells.Item(20, 4).Value = -14
$ws.Cells.Item(20, 5).Value = -7.2
$ws.Cells.Item(20, 6).Value = 17.73
$ws.Cells.Item(21, 1).Value = "RM 135"
$ws.Cells.Item(21, 2).Value = -18.9
$ws.Cells.Item(21, 3).Value = 12.7
$ws.Cells.Item(21, 4).Value = -14.3
$ws.Cells.Item(21, 5).Value = -8.7
$ws.Cells.Item(21, 6).Value = 16.58
$ws.Cells.Item(22, 1).Value = "RM 138"
$ws.Cells.Item(22, 2).Value = -19.3
$ws.Cells.Item(22, 3).Value = ""
$ws.Cells.Item(22, 4).Value = ""
$ws.Cells.Item(22, 5).Value = -6.1
$ws.Cells.Item(22, 6).Value = 16.81
$ws.Cells.Item(23, 1).Value = "RM 140"
$ws.Cells.Item(23, 2).Value = -19.5
$ws.Cells.Item(23, 3).Value = ""
$ws.Cells.Item(23, 4).Value = -13.9
$ws.Cells.Item(23, 5).Value = -7
$ws.Cells.Item(23, 6).Value = 16.48
$ws.Cells.Item(24, 1).Value = "RM 142a"
$ws.Cells.Item(24, 2).Value = -17.7
$ws.Cells.Item(24, 3).Value = ""
$ws.Cells.Item(24, 4).Value = -13.9
$ws.Cells.Item(24, 5).Value = -8.1
$ws.Cells.Item(24, 6).Value = 16.78
$ws.Cells.Item(25, 1).Value = "RM 145"
$ws.Cells.Item(25, 2).Value = -19.5
$ws.Cells.Item(25, 3).Value = 10.7
$ws.Cells.Item(25, 4).Value = -15.5
$ws.Cells.Item(25, 5).Value = -7.1
$ws.Cells.Item(25, 6).Value = ""
$ws.Cells.Item(26, 1).Value = "SC 5"
$ws.Cells.Item(26, 2).Value = -20.2
$ws.Cells.Item(26, 3).Value = 10.8
$ws.Cells.Item(26, 4).Value = -13.8
$ws.Cells.Item(26, 5).Value = -5
$ws.Cells.Item(26, 6).Value = 17.38
$ws.Cells.Item(27, 1).Value = "SC 101"
$ws.Cells.Item(27, 2).Value = ""
$ws.Cells.Item(27, 3).Value = 10
$ws.Cells.Item(27, 4).Value = ""
$ws.Cells.Item(27, 5).Value = -10
$ws.Cells.Item(27, 6).Value = 17
$ws.Cells.Item(28, 1).Value = "SC 105"
$ws.Cells.Item(28, 2).Value = ""
$ws.Cells.Item(28, 3).Value = 11.1
$ws.Cells.Item(28, 4).Value = -13.7
$ws.Cells.Item(28, 5).Value = -5.9
$ws.Cells.Item(28, 6).Value = 17.44
$ws.Cells.Item(29, 1).Value = "SC 119"
$ws.Cells.Item(29, 2).Value = -19.5
$ws.Cells.Item(29, 3).Value = 11.2
$ws.Cells.Item(29, 4).Value = ""
$ws.Cells.Item(29, 5).Value = -6.8
$ws.Cells.Item(29, 6).Value = 18.06
$ws.Cells.Item(30, 1).Value = "SC 120"
$ws.Cells.Item(30, 2).Value = -19.7
$ws.Cells.Item(30, 3).Value = 11.4
$ws.Cells.Item(30, 4).Value = -13.6
$ws.Cells.Item(30, 5).Value = -5.7
$ws.Cells.Item(30, 6).Value = ""
$ws.Cells.Item(31, 1).Value = "SC 132"
$ws.Cells.Item(31, 2).Value = ""
$ws.Cells.Item(31, 3).Value = 15.3
$ws.Cells.Item(31, 4).Value = -13.7
$ws.Cells.Item(31, 5).Value = ""
$ws.Cells.Item(31, 6).Value = 17.18
$ws.Cells.Item(32, 1).Value = "SC 193"
$ws.Cells.Item(32, 2).Value = ""
$ws.Cells.Item(32, 3).Value = 10.5
$ws.Cells.Item(32, 4).Value = -14.7
$ws.Cells.Item(32, 5).Value = -6.4
$ws.Cells.Item(32, 6).Value = 17.39
$ws.Cells.Item(33, 1).Value = "SC 232"
$ws.Cells.Item(33, 2).Value = -19.5
$ws.Cells.Item(33, 3).Value = 10.4
$ws.Cells.Item(33, 4).Value = -14.1
$ws.Cells.Item(33, 5).Value = -10.7
$ws.Cells.Item(33, 6).Value = 17.53
